# Apply the edits described by the diff:
# 1. Update the absPath url in the workbook (add a space in "ПоказателиЦУР" -> "Показатели ЦУР")
# 2. On sheet1:
#    - extend dimension/used range to include column Q
#    - change sheet view: topLeftCell=C1, selection Q12
#    - add a new "2020" column (Q) with data for rows 3-8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the absolute path stored in the workbook metadata ---
$wb.AbsolutePathUrl = "C:\Users\korozbaeva\Desktop\Показатели ЦУР для Платформы\Глобальные показатели ЦУР\"

# --- 2. Add the new 2020 column (Q) with values/formulas ---
$ws.Range("Q3").Value = 2020
$ws.Range("Q3").Style = $ws.Range("P3").Style

$ws.Range("Q6").Value = 312
$ws.Range("Q6").Style = $ws.Range("P6").Style

$ws.Range("Q7").Value = 1856
$ws.Range("Q7").Style = $ws.Range("P7").Style

$ws.Range("Q8").Value = 4337617
$ws.Range("Q8").Style = $ws.Range("P8").Style

$ws.Range("Q4").Formula = "=Q6/Q8*100000"
$ws.Range("Q4").Style = $ws.Range("P4").Style

$ws.Range("Q5").Formula = "=Q7/Q8*100000"
$ws.Range("Q5").Style = $ws.Range("P5").Style

# --- 3. Update the sheet view: scroll to column C and select Q12 ---
$ws.Range("Q12").Select()
$excel.ActiveWindow.ScrollColumn = 3

$wb.Save()
